# The order sheet got a batch of "Natalie's" juice items that were previously
# dropped onto the wrong sheet/layer. Append them as new rows 14-18, matching
# the existing layout (SKU, Name, Quantity, Cost Per, Total Cost) where every
# value - including the numeric-looking ones - is stored as plain text, same
# as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("TN374", "Natalie's - Lemonade",           "1", "9.30",  "9.30"),
    @("AH252", "Natalie's - Orange Juice",       "1", "24.50", "24.50"),
    @("TN454", "Natalie's - Orange Mango",       "1", "13.38", "13.38"),
    @("TN362", "Natalie's - Orange Pineapple",   "1", "13.38", "13.38"),
    @("TN380", "Natalie's - Strawberry Lemonade","1", "10.15", "10.15")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Quantity/Cost Per/Total Cost look numeric ("1", "9.30", ...) but the
    # sheet stores every column as text, so force text format on those
    # columns first - otherwise Excel would coerce them into numbers and
    # drop the trailing zero (9.30 -> 9.3).
    $ws.Range("C$r`:E$r").NumberFormat = "@"

    $ws.Range("A$r").Value = $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
    $ws.Range("E$r").Value = $values[4]
}
